$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(33, 1).Value = 616
$ws.Cells.Item(33, 2).Value = 100
$ws.Cells.Item(33, 3).Value = 100
$ws.Cells.Item(33, 4).Value = 100
$ws.Cells.Item(33, 5).Value = "Wood"
